$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the "Cases" query text (row 2, column B / "query") ---
# Removed the trailing "Cohort" return line and the bug now truncates the
# "Response to Treatment" column alias (matches the author's commit).
$newCasesQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
WHERE demo.breed IN ['American Staffordshire Terrier']
MATCH (c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment
'@

# --- Fix the "Files" query text (row 4, column B / "query") ---
# Swapped the last return column from Study Code to Cohort.
$newFilesQuery = @'
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
WHERE demo.breed IN ['American Staffordshire Terrier']
OPTIONAL MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
WITH DISTINCT f, parent, c, demo, diag, s
RETURN coalesce(f.file_name, '') AS `File Name`, 
        coalesce(f.file_type, '') AS `File Type`, 
        coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
        coalesce(f.file_format, '') AS `File Format`,
        coalesce(f.file_size, '') AS `Size`,
        coalesce(c.case_id, '') AS `Case ID`, 
        coalesce(demo.breed,'') AS Breed , 
        coalesce(diag.disease_term,'') AS Diagnosis,
        coalesce(co.cohort_description, '') AS `Cohort`
'@

# The here-string adds a trailing newline after the last line; strip it so
# the cell text matches the source exactly (no trailing newline).
$newCasesQuery = $newCasesQuery.TrimEnd("`r","`n")
$newFilesQuery = $newFilesQuery.TrimEnd("`r","`n")

$ws.Cells.Item(2, 2).Value = $newCasesQuery
$ws.Cells.Item(4, 2).Value = $newFilesQuery

# --- View: zoom in from 70% to 145% ---
$excel.ActiveWindow.Zoom = 145

# --- Row heights updated to reflect the new wrapped-text line counts ---
$ws.Rows.Item(2).RowHeight = 244.8
$ws.Rows.Item(3).RowHeight = 216
$ws.Rows.Item(4).RowHeight = 244.8
